$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update selection (B23 -> B22)
$ws.Range("B22").Select()

# B2 text stays the same ("OpenTBS demo"); shared string index changes automatically.
$ws.Range("B2").Value = "OpenTBS demo"

# B13: "Merging Microsoft Excel templates..." text updated
$ws.Range("B13").Value = "Merging Microsoft Excel templates with OpenTBS has several big limitations because of the OpenXML format for Excel."

# B14: new text about Excel formulas
$ws.Range("B14").Value = "* Excel formulas are saved in Sheet subfiles with both the expression and the result. This means that when you write one TBS tag in Excel, it may have two in the XML source of the template."

# B15: Texts in cells note, updated wording
$ws.Range("B15").Value = "* Texts in cells are not saved in the sheet subfile but in the 'xl/sharedStrings.xml' subfile. This means you cannot refer to any sheet tag (like row) for TBS tags placed in a text cell."

# B17: Pictures placed note, updated wording
$ws.Range("B17").Value = "* Pictures placed in the sheet are not referenced in the sheet subfile but in another XML subfile (for instance 'xl\drawings\drawing1.xml'). This means you cannot use the usual parameter ""ope=changepic"" to change pictures in a sheet."

# B18: Cells are saved note, updated wording
$ws.Range("B18").Value = "* Cells are saved in sheets with their absolute position (row+column). Thus, when you merge them using MergeBlock(), positions are duplicated. This can produce wrong sheets."

# B20: Far all those raisons note, updated wording
$ws.Range("B20").Value = "Far all those raisons, its seems than it is not possible to use MergeBlock() an Excel template."

# B21 text stays the same ("You can anly use simple TBS fields...")
$ws.Range("B21").Value = "You can anly use simple TBS fields. Which is quite few."
